$d = $word.ActiveDocument

# Remove the "ResponseId" (Heading2) paragraph and the "Hello" paragraph
# that follow the title, leaving only "Logan's Report".
$start = $d.Paragraphs.Item(2).Range.Start
$end = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
$rng = $d.Range($start, $end)
$rng.Delete()
